$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "14 - Diplomacy" protocol block (columns K = Peace, L = War)
# Row 4 - Parametres: add the new diplomacy protocol parametres
$ws.Range("K4").Value = "country, country"
$ws.Range("L4").Value = "country, country"

# L4 picks up K4's border/format (matches the corrected layout in the diff)
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 5 - Returns: corrected country death protocol
$ws.Range("K5").Value = "-"
$ws.Range("L5").Value = "-"

# Row 6 - Other
$ws.Range("K6").Value = "-"
$ws.Range("L6").Value = "-"

# Move / update the active selection to A6
$ws.Range("A6").Select()
